# Auto-generated edit script: updates odds/score-count cells per the commit diff
# for Jogos_da_Semana_FlashScore_2025-01-29.xlsx (Sheet1).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("Q3").Value = 1.99
$ws.Range("R3").Value = 1.74

# Row 4
$ws.Range("U4").Value = 5.2

# Row 7
$ws.Range("O7").Value = 1.2
$ws.Range("P7").Value = 4.5
$ws.Range("W7").Value = 2.63
$ws.Range("X7").Value = 1.5

# Row 8
$ws.Range("G8").Value = 2.15
$ws.Range("H8").Value = 3.2
$ws.Range("I8").Value = 3.1
$ws.Range("J8").Value = 2.88
$ws.Range("S8").Value = 1.88
$ws.Range("T8").Value = 1.98
$ws.Range("AF8").Value = 21
$ws.Range("AH8").Value = 26

# Row 9
$ws.Range("N9").Value = 15
$ws.Range("T9").Value = 2.3
$ws.Range("U9").Value = 1.98
$ws.Range("V9").Value = 1.83
$ws.Range("W9").Value = 2.5
$ws.Range("X9").Value = 1.5

# Row 10
$ws.Range("G10").Value = 2.75
$ws.Range("S10").Value = 2.1
$ws.Range("T10").Value = 1.67
$ws.Range("AG10").Value = 23
$ws.Range("AM10").Value = 8
$ws.Range("AN10").Value = 12
$ws.Range("AQ10").Value = 21
$ws.Range("AS10").Value = 301

# Row 11
$ws.Range("G11").Value = 1.5
$ws.Range("H11").Value = 4.75
$ws.Range("J11").Value = 2
$ws.Range("K11").Value = 2.6
$ws.Range("M11").Value = 1.02
$ws.Range("N11").Value = 19
$ws.Range("O11").Value = 1.14
$ws.Range("P11").Value = 5.5
$ws.Range("S11").Value = 1.47
$ws.Range("T11").Value = 2.5
$ws.Range("U11").Value = 1.83
$ws.Range("V11").Value = 2.03
$ws.Range("W11").Value = 2.25
$ws.Range("X11").Value = 1.57
$ws.Range("Y11").Value = 1.25
$ws.Range("Z11").Value = 3.75
$ws.Range("AA11").Value = 1.67
$ws.Range("AB11").Value = 2.1
$ws.Range("AC11").Value = 9.5
$ws.Range("AD11").Value = 8.5
$ws.Range("AG11").Value = 11
$ws.Range("AI11").Value = 19
$ws.Range("AJ11").Value = 10
$ws.Range("AK11").Value = 17
$ws.Range("AM11").Value = 19
$ws.Range("AR11").Value = 34
$ws.Range("AS11").Value = 151

# Row 12
$ws.Range("G12").Value = 3.5
$ws.Range("H12").Value = 3.2
$ws.Range("K12").Value = 2.1
$ws.Range("L12").Value = 2.88
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 9.5
$ws.Range("O12").Value = 1.3
$ws.Range("P12").Value = 3.4
$ws.Range("S12").Value = 2.05
$ws.Range("T12").Value = 1.8
$ws.Range("W12").Value = 3.5
$ws.Range("X12").Value = 1.29
$ws.Range("Y12").Value = 1.44
$ws.Range("Z12").Value = 2.63
$ws.Range("AA12").Value = 1.8
$ws.Range("AB12").Value = 1.95
$ws.Range("AC12").Value = 10
$ws.Range("AF12").Value = 34
$ws.Range("AH12").Value = 34
$ws.Range("AI12").Value = 9.5
$ws.Range("AK12").Value = 13
$ws.Range("AM12").Value = 7.5
$ws.Range("AO12").Value = 9
$ws.Range("AQ12").Value = 17
$ws.Range("AR12").Value = 29
$ws.Range("AS12").Value = 251

# Row 13
$ws.Range("T13").Value = 1.75

# Row 14
$ws.Range("O14").Value = 1.29
$ws.Range("P14").Value = 3.5
$ws.Range("S14").Value = 1.93
$ws.Range("T14").Value = 1.93

# Row 16
$ws.Range("Q16").Value = 1.85
$ws.Range("R16").Value = 1.95
$ws.Range("S16").Value = 2.5
$ws.Range("T16").Value = 1.5
$ws.Range("W16").Value = 5
$ws.Range("X16").Value = 1.17

# Row 17
$ws.Range("G17").Value = 2.5
$ws.Range("H17").Value = 3.65
$ws.Range("I17").Value = 2.45
$ws.Range("K17").Value = 2.4
$ws.Range("L17").Value = 2.85
$ws.Range("P17").Value = 5.2
$ws.Range("S17").Value = 1.4
$ws.Range("T17").Value = 2.72
$ws.Range("W17").Value = 1.95
$ws.Range("X17").Value = 1.75
$ws.Range("Y17").Value = 1.24
$ws.Range("AA17").Value = 1.36
$ws.Range("AB17").Value = 2.9
$ws.Range("AJ17").Value = 8.25
$ws.Range("AK17").Value = 10.25
$ws.Range("AM17").Value = 15.5
$ws.Range("AN17").Value = 18.5
$ws.Range("AO17").Value = 10

# Row 18
$ws.Range("G18").Value = 4.55
$ws.Range("J18").Value = 4.55
$ws.Range("K18").Value = 2.35
$ws.Range("N18").Value = 9
$ws.Range("O18").Value = 1.18
$ws.Range("P18").Value = 4.3
$ws.Range("S18").Value = 1.55
$ws.Range("T18").Value = 2.3
$ws.Range("W18").Value = 2.35
$ws.Range("X18").Value = 1.53
$ws.Range("Y18").Value = 1.3
$ws.Range("Z18").Value = 3.25
$ws.Range("AA18").Value = 1.6
$ws.Range("AB18").Value = 2.22
$ws.Range("AC18").Value = 17.5
$ws.Range("AD18").Value = 30
$ws.Range("AE18").Value = 14.5
$ws.Range("AF18").Value = 80
$ws.Range("AG18").Value = 37
$ws.Range("AH18").Value = 35
$ws.Range("AI18").Value = 9
$ws.Range("AJ18").Value = 8
$ws.Range("AL18").Value = 45
$ws.Range("AM18").Value = 9
$ws.Range("AN18").Value = 9
$ws.Range("AO18").Value = 8
$ws.Range("AP18").Value = 13.5
$ws.Range("AR18").Value = 19.5

# Row 19
$ws.Range("H19").Value = 2.88
$ws.Range("I19").Value = 2.88
$ws.Range("M19").Value = 1.11
$ws.Range("N19").Value = 6.5
$ws.Range("AF19").Value = 26
$ws.Range("AJ19").Value = 5.5
$ws.Range("AM19").Value = 7.5
$ws.Range("AN19").Value = 13

# Row 20
$ws.Range("G20").Value = 5.4
$ws.Range("H20").Value = 4.25
$ws.Range("I20").Value = 1.52
$ws.Range("J20").Value = 5.2
$ws.Range("K20").Value = 2.35
$ws.Range("L20").Value = 2.02
$ws.Range("N20").Value = 8.5
$ws.Range("P20").Value = 3.9
$ws.Range("X20").Value = 1.45
$ws.Range("Y20").Value = 1.33
$ws.Range("Z20").Value = 3.05
$ws.Range("AA20").Value = 1.78
$ws.Range("AB20").Value = 1.93
$ws.Range("AC20").Value = 16.5
$ws.Range("AD20").Value = 32
$ws.Range("AF20").Value = 100
$ws.Range("AH20").Value = 50
$ws.Range("AI20").Value = 8.5
$ws.Range("AJ20").Value = 8.25
$ws.Range("AK20").Value = 16.5
$ws.Range("AL20").Value = 70
$ws.Range("AM20").Value = 7.7
$ws.Range("AO20").Value = 8.25
$ws.Range("AP20").Value = 10.75
$ws.Range("AQ20").Value = 11.5
$ws.Range("AR20").Value = 24
$ws.Range("AS20").Value = 500
